$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3; this shifts the former rows 3-9 down to 4-10.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with a new weekly price observation for Kiwi.
# Columns A,B,C,E,F,G,H,I,J,K,L,R stay constant across the whole subset,
# so copy them from row 2 (same "Agrícola del Norte S.A. de Arica" / Kiwi series).
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = 44602
$ws.Range("E3").Value = $ws.Range("E2").Value()
$ws.Range("F3").Value = $ws.Range("F2").Value()
$ws.Range("G3").Value = $ws.Range("G2").Value()
$ws.Range("H3").Value = $ws.Range("H2").Value()
$ws.Range("I3").Value = $ws.Range("I2").Value()
$ws.Range("J3").Value = $ws.Range("J2").Value()
$ws.Range("K3").Value = $ws.Range("K2").Value()
$ws.Range("L3").Value = $ws.Range("L2").Value()
$ws.Range("M3").Value = 270
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("Q3").Value = "$/bandeja 18 kilos"
$ws.Range("R3").Value = $ws.Range("R2").Value()
$ws.Range("S3").Value = 1139
$ws.Range("T3").Value = 18
